$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.830.38'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '3.139.88'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '574.43'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").Value = '164.75'
$ws.Range("E6").Value = '  -2.67%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '0.575'
$ws.Range("E8").Value = '  -5.19%  '
$ws.Range("D9").Value = '3.154.02'
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("D10").Value = '0.118'
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("E11").Value = '  -2.68%  '
$ws.Range("D12").Value = '0.382'
$ws.Range("E12").Value = '  -2.29%  '
$ws.Range("D13").Value = '3.686.39'
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("E14").Value = '  -2.03%  '
$ws.Range("D15").Value = '64.842.01'
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").Value = '25.02'
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("D17").Value = '3.144.76'
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").Value = '0.0000156'
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("D19").Value = '414.21'
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("E20").Value = '  -1.75%  '
$ws.Range("D21").Value = '12.49'
$ws.Range("E21").Value = '  -3.64%  '
$ws.Range("D22").Value = '7.04'
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '68.84'
$ws.Range("E24").Value = '  -1.88%  '
$ws.Range("D25").Value = '0.483'
$ws.Range("E25").Value = '  -3.00%  '
$ws.Range("D26").Value = '0.193'
$ws.Range("E26").Value = '  -5.20%  '
$ws.Range("D27").Value = '0.0000104'
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").Value = '9.10'
$ws.Range("E28").Value = '  +3.28%  '
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("D32").Value = '21.31'
$ws.Range("E32").Value = '  -2.09%  '
$ws.Range("D33").Value = '163.28'
$ws.Range("E33").Value = '  +4.85%  '
$ws.Range("D34").Value = '4.88'
$ws.Range("E34").Value = '  -3.64%  '
$ws.Range("D35").Value = '6.25'
$ws.Range("E35").Value = '  -1.93%  '
$ws.Range("D36").Value = '1.13'
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").Value = '1.36'
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("D38").Value = '1.69'
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("D39").Value = '2.623.61'
$ws.Range("E39").Value = '  -3.21%  '
$ws.Range("D40").Value = '4.15'
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("D41").Value = '23.82'
$ws.Range("E41").Value = '  -1.83%  '
$ws.Range("D42").Value = '38.38'
$ws.Range("E42").Value = '  -1.83%  '
$ws.Range("D43").Value = '0.693'
$ws.Range("E43").Value = '  -3.39%  '
$ws.Range("D44").Value = '0.0616'
$ws.Range("E44").Value = '  -1.12%  '
$ws.Range("D45").Value = '5.32'
$ws.Range("E45").Value = '  -3.68%  '
$ws.Range("D46").Value = '291.47'
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("D47").Value = '0.0256'
$ws.Range("E47").Value = '  -3.06%  '
$ws.Range("D48").Value = '21.29'
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("D49").Value = '0.996'
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").Value = '0.0977'
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '1.92'
$ws.Range("E51").Value = '  -4.13%  '
